# Generate Report for Handoff
#
# The "b.md" source file has now been handed off again (a newer source
# changed since the last handback), so update the "b.md" row on the
# Overview sheet and on each per-locale (zh-cn / de-de) detail sheet:
#   - Status moves from "Handed back: in sync with en-US" to
#     "Ready for handoff"
#   - Content Duplicate flips from True to False
#   - Latest Handoff File / Datetime point at the new handoff package
#   - Error Detail now explains the handback file is stale

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$text) {
    # Force a plain text cell (t="s") instead of letting Excel infer a
    # Boolean/number type for tokens like "True"/"False". The leading
    # apostrophe marks the entry as text (quote-prefixed); resetting the
    # style back to Normal afterwards drops the quote-prefix flag again so
    # the cell ends up indistinguishable from an ordinarily authored
    # string cell.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# --- Overview sheet: row 3 is "b.md" ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-09-04 02:41:25"

# --- zh-cn sheet: row 3 is "b.md" ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Ready for handoff"
Set-TextValue $ws.Range("F3") "False"
$ws.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$ws.Range("H3").Value = "2016-09-04 02:41:20"
$ws.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/816db466cb003b699058d14ac759f5769d88a66e/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/31b50ccec0ae80bd393341f0b99e12674f1d8470/e2e/b.md."
# Widen the now-populated "Error Detail" column to a round 40 characters
# (ColumnWidth includes ~0.8333 of cell-padding overhead on top of the
# character count that ends up in the saved <col width="..."> value).
$ws.Columns.Item(16).ColumnWidth = 39.1666666666667

# --- de-de sheet: row 3 is "b.md" ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "Ready for handoff"
Set-TextValue $ws.Range("F3") "False"
$ws.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$ws.Range("H3").Value = "2016-09-04 02:41:25"
$ws.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/816db466cb003b699058d14ac759f5769d88a66e/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/31b50ccec0ae80bd393341f0b99e12674f1d8470/e2e/b.md."
$ws.Columns.Item(16).ColumnWidth = 39.1666666666667
